# Auto-generated edit script: updates crafting-profit calculation outputs
# (columns H-N) across multiple sheets, per scheduled profit-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8028.9165
$ws.Range("J62").Value = 9281.071
$ws.Range("L62").Value = 9281.071
$ws.Range("N62").Value = -10529.071
$ws.Range("H65").Value = 8028.9165
$ws.Range("J65").Value = 9281.071
$ws.Range("L65").Value = 46405.355
$ws.Range("N65").Value = -52645.355
$ws.Range("H88").Value = 1755.48
$ws.Range("I88").Value = 2111
$ws.Range("K88").Value = 2111
$ws.Range("M88").Value = -1705
$ws.Range("H91").Value = 1755.48
$ws.Range("I91").Value = 2111
$ws.Range("K91").Value = 2111
$ws.Range("M91").Value = -707
$ws.Range("H138").Value = 2804.9412
$ws.Range("I138").Value = 2329.125
$ws.Range("J138").Value = 3227.889
$ws.Range("K138").Value = 6987.375
$ws.Range("L138").Value = 9683.667000000001
$ws.Range("M138").Value = -1847.375
$ws.Range("N138").Value = -19963.667
$ws.Range("H141").Value = 5780
$ws.Range("I141").Value = 6656.1
$ws.Range("K141").Value = 19968.3
$ws.Range("M141").Value = -14788.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3896
$ws.Range("I61").Value = 3896
$ws.Range("K61").Value = 3896
$ws.Range("M61").Value = -3684
$ws.Range("H102").Value = 4399.643
$ws.Range("I102").Value = 3739.4
$ws.Range("J102").Value = 6050.25
$ws.Range("K102").Value = 3739.4
$ws.Range("L102").Value = 6050.25
$ws.Range("M102").Value = -2117.4
$ws.Range("N102").Value = -9294.25
$ws.Range("H132").Value = 25405.38
$ws.Range("I132").Value = 25962.44
$ws.Range("J132").Value = 2566
$ws.Range("K132").Value = 77887.31999999999
$ws.Range("L132").Value = 7698
$ws.Range("M132").Value = -75357.31999999999
$ws.Range("N132").Value = -12758
$ws.Range("H136").Value = 3896
$ws.Range("I136").Value = 3896
$ws.Range("K136").Value = 11688
$ws.Range("M136").Value = -9138

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 100764
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H94").Value = 992
$ws.Range("I94").Value = 859.5925999999999
$ws.Range("J94").Value = 2183.6667
$ws.Range("K94").Value = 859.5925999999999
$ws.Range("L94").Value = 2183.6667
$ws.Range("M94").Value = -408.5925999999999
$ws.Range("N94").Value = -3085.6667
$ws.Range("H122").Value = 150000
$ws.Range("J122").Value = 150000
$ws.Range("L122").Value = 150000
$ws.Range("N122").Value = -159800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2582.5652
$ws.Range("I31").Value = 2055.4375
$ws.Range("J31").Value = 3787.4285
$ws.Range("K31").Value = 2055.4375
$ws.Range("L31").Value = 3787.4285
$ws.Range("M31").Value = -1760.4375
$ws.Range("N31").Value = -4377.4285
$ws.Range("H34").Value = 2582.5652
$ws.Range("I34").Value = 2055.4375
$ws.Range("J34").Value = 3787.4285
$ws.Range("K34").Value = 2055.4375
$ws.Range("L34").Value = 3787.4285
$ws.Range("M34").Value = -1853.4375
$ws.Range("N34").Value = -4191.4285
$ws.Range("H60").Value = 35000
$ws.Range("I60").Value = 35000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 35000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -34489
$ws.Range("N60").ClearContents()
$ws.Range("H132").Value = 1190.8206
$ws.Range("I132").Value = 1201.1052
$ws.Range("K132").Value = 3603.3156
$ws.Range("M132").Value = -1073.3156
$ws.Range("H134").Value = 60315.89
$ws.Range("I134").Value = 86989.586
$ws.Range("K134").Value = 260968.758
$ws.Range("M134").Value = -258433.758

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2166.5
$ws.Range("I109").Value = 1999.8
$ws.Range("K109").Value = 5999.4
$ws.Range("M109").Value = -4959.4
$ws.Range("H131").Value = 1728499.1
$ws.Range("I131").Value = 1786.2727
$ws.Range("J131").Value = 2132623.2
$ws.Range("K131").Value = 5358.8181
$ws.Range("L131").Value = 6397869.600000001
$ws.Range("M131").Value = -318.8181000000004
$ws.Range("N131").Value = -6407949.600000001
$ws.Range("H132").Value = 1133
$ws.Range("I132").Value = 1133
$ws.Range("K132").Value = 10197
$ws.Range("M132").Value = -7667
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H140").Value = 3250.6924
$ws.Range("I140").Value = 2846.2
$ws.Range("K140").Value = 8538.599999999999
$ws.Range("M140").Value = -3358.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 71221.60000000001
$ws.Range("J123").Value = 71221.60000000001
$ws.Range("L123").Value = 71221.60000000001
$ws.Range("N123").Value = -76121.60000000001
$ws.Range("H126").Value = 10069.6
$ws.Range("I126").Value = 8349.286
$ws.Range("J126").Value = 11574.875
$ws.Range("K126").Value = 25047.858
$ws.Range("L126").Value = 34724.625
$ws.Range("M126").Value = -22577.858
$ws.Range("N126").Value = -39664.625
$ws.Range("H132").Value = 35913.62
$ws.Range("I132").Value = 39861.348
$ws.Range("J132").Value = 1700
$ws.Range("K132").Value = 119584.044
$ws.Range("L132").Value = 5100
$ws.Range("M132").Value = -117054.044
$ws.Range("N132").Value = -10160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 9697
$ws.Range("J97").Value = 9697
$ws.Range("L97").Value = 9697
$ws.Range("N97").Value = -11679
$ws.Range("H132").Value = 123827.7
$ws.Range("I132").Value = 153360.5
$ws.Range("K132").Value = 460081.5
$ws.Range("M132").Value = -457551.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H96").Value = 1860.6
$ws.Range("I96").Value = 2075.75
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 2075.75
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = -702.75
$ws.Range("N96").Value = -3746
$ws.Range("H122").Value = 1307.8334
$ws.Range("J122").Value = 1600
$ws.Range("L122").Value = 4800
$ws.Range("N122").Value = -9700
$ws.Range("H132").Value = 28546.83
$ws.Range("I132").Value = 34669.793
$ws.Range("K132").Value = 104009.379
$ws.Range("M132").Value = -101479.379
$ws.Range("H135").Value = 68857.5
$ws.Range("J135").Value = 68857.5
$ws.Range("L135").Value = 68857.5
$ws.Range("N135").Value = -78997.5
